$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert the filing/publication date columns (D & E) from date-serial
#     cells to plain text cells formatted "dd/mm/yyyy" -----------------------
# Apply a text number format first so the subsequent string assignments are
# not re-interpreted as date serials. Column D (including its header) carries
# the text style; column E only needs it on the data rows.
$ws.Range("D1:D14").NumberFormat = "@"
$ws.Range("E2:E14").NumberFormat = "@"

$dates = @{
    2  = @("11/03/2013", "13/06/2013")
    3  = @("14/05/2013", "15/08/2013")
    4  = @("17/07/2013", "13/09/2013")
    5  = @("19/09/2013", "15/12/2013")
    6  = @("22/11/2013", "23/02/2014")
    7  = @("25/01/2014", "17/04/2014")
    8  = @("30/03/2014", "13/06/2014")
    9  = @("02/06/2014", "05/08/2014")
    10 = @("05/08/2014", "13/11/2014")
    11 = @("08/10/2014", "09/01/2015")
    12 = @("11/12/2014", "24/03/2015")
    13 = @("13/02/2015", "15/06/2015")
    14 = @("18/04/2015", "13/06/2015")
}

for ($row = 2; $row -le 14; $row++) {
    $vals = $dates[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}

# --- Drop the trailing blank rows 15-23 (only D15:D23 held leftover styling)
$ws.Rows("15:23").Delete()

# --- Move the active selection to E15, matching the saved view state -------
$ws.Range("E15").Select()
